# Updated cryptos list on Mon Apr 24 15:25:18 UTC 2023 with GitHub Actions
#
# For every "Price" (column D) cell, the new value looks numeric to Excel's
# text-parser (single dot decimals, multi-dot "thousand.thousand.cents"
# figures, etc.), so a plain .Value assignment would silently convert it to
# a real number (and round/reformat it). To keep these cells as literal text
# - matching the source workbook's inlineStr cells - each D-column write is
# wrapped: force Text format, assign the literal string, then reset the
# cell's style back to "Normal" so no stray number-format style lingers on
# the cell (only the cell's value changes, nothing structural).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.642.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.015"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.014"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4659"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3929"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07975"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.003"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.872.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.963"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.227"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.016"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06729"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001044"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.010"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.677.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.468"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.304"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.097.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.138"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.443"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9771"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09453"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.626"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.306"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.341"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06043"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02235"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.200"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.318"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.012"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5953"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1874"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.251"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5632"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.928"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06754"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.054"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.43%  "
